# Generate Report for Handback
#
# Fills in the "Latest Target File" (I) and "Latest Handback File" (J) columns
# for both localized-language tabs (zh-cn, de-de), links the new "Latest Target
# File" cells back to the source markdown (mirroring the existing column A
# hyperlinks), stamps the "Latest Handback DateTime" (K) with the handback
# timestamp (distinct per language), and flips the shared "Status" text from
# "Ready for handoff" to "Handed back: in sync with en-US" everywhere it
# appears (Overview zh-cn/de-de columns + the Status column on each language
# tab). Column widths on the touched columns are widened to fit the new text.

$wb = $excel.ActiveWorkbook

# Cornflower blue (FF6495ED), matching the existing hyperlink font exactly.
$hyperlinkColor = 15570276
$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef4c4f86babbf720b66eec4f1b1dec20daac7375/e2e/"

# ---------------------------------------------------------------------------
# zh-cn tab
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$cell = $wsZh.Range("I2")
$wsZh.Hyperlinks.Add($cell, ($repoBase + "71107d88-d211-4f7c-839e-e0f6ec9a8aa0.md"), "", "", "71107d88-d211-4f7c-839e-e0f6ec9a8aa0.md")
$cell.Font.Underline = $true
$cell.Font.Color = $hyperlinkColor
$wsZh.Range("J2").Value = "71107d88-d211-4f7c-839e-e0f6ec9a8aa0.e06e07c1042e87b2bf3e745696dce36450f2941d.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-30 18:59:48"

$cell = $wsZh.Range("I3")
$wsZh.Hyperlinks.Add($cell, ($repoBase + "af3f0b90-2ce4-4fe3-aa15-a8a3920e9e6d.md"), "", "", "af3f0b90-2ce4-4fe3-aa15-a8a3920e9e6d.md")
$cell.Font.Underline = $true
$cell.Font.Color = $hyperlinkColor
$wsZh.Range("J3").Value = "af3f0b90-2ce4-4fe3-aa15-a8a3920e9e6d.5bcf0b977c20f05706fe307b0ca38c854d1b72e4.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-30 18:59:48"

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Columns.Item(3).ColumnWidth = 29.15
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de tab
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$cell = $wsDe.Range("I2")
$wsDe.Hyperlinks.Add($cell, ($repoBase + "71107d88-d211-4f7c-839e-e0f6ec9a8aa0.md"), "", "", "71107d88-d211-4f7c-839e-e0f6ec9a8aa0.md")
$cell.Font.Underline = $true
$cell.Font.Color = $hyperlinkColor
$wsDe.Range("J2").Value = "71107d88-d211-4f7c-839e-e0f6ec9a8aa0.e06e07c1042e87b2bf3e745696dce36450f2941d.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-30 18:59:55"

$cell = $wsDe.Range("I3")
$wsDe.Hyperlinks.Add($cell, ($repoBase + "af3f0b90-2ce4-4fe3-aa15-a8a3920e9e6d.md"), "", "", "af3f0b90-2ce4-4fe3-aa15-a8a3920e9e6d.md")
$cell.Font.Underline = $true
$cell.Font.Color = $hyperlinkColor
$wsDe.Range("J3").Value = "af3f0b90-2ce4-4fe3-aa15-a8a3920e9e6d.5bcf0b977c20f05706fe307b0ca38c854d1b72e4.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-30 18:59:55"

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Columns.Item(3).ColumnWidth = 29.15
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# Overview tab - same "Status" text swap for the zh-cn / de-de columns
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

Write-Output "Handback report generated."
